$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002837
$ws.Range("H2").Value = 0.008510999999999999
$ws.Range("I2").Value = 0.000071082477304929285628898317
$ws.Range("J2").Value = 0.000071082477304929285628898317
$ws.Range("M2").Value = 1.815761
$ws.Range("N2").Value = 5.447283000000001
$ws.Range("O2").Value = 0.07007596730428067
$ws.Range("P2").Value = 0.07007596730428067
$ws.Range("Q2").Value = 0.005151313957
$ws.Range("R2").Value = 0.046361825613
$ws.Range("S2").Value = 0.000004981173355527497247888409
$ws.Range("T2").Value = 0.000004981173355527497247888409

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002837
$ws.Range("H3").Value = 0.008510999999999999
$ws.Range("I3").Value = 0.000071082477304929285628898317
$ws.Range("J3").Value = 0.000071082477304929285628898317
$ws.Range("O3").Value = 0.5079540516959071
$ws.Range("P3").Value = 0.5079540516959072
$ws.Range("Q3").Value = 0.03733991690266666
$ws.Range("R3").Value = 0.336059252124
$ws.Range("S3").Value = 0.000036106632351621203123881454
$ws.Range("T3").Value = 0.000036106632351621203123881454

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002837
$ws.Range("H4").Value = 0.008510999999999999
$ws.Range("I4").Value = 0.000071082477304929285628898317
$ws.Range("J4").Value = 0.000071082477304929285628898317
$ws.Range("M4").Value = 9.711409333333334
$ws.Range("N4").Value = 29.134228
$ws.Range("O4").Value = 0.3747940411327002
$ws.Range("P4").Value = 0.3747940411327002
$ws.Range("Q4").Value = 0.02755126827866667
$ws.Range("R4").Value = 0.247961414508
$ws.Range("S4").Value = 0.000026641288922837888845061849
$ws.Range("T4").Value = 0.000026641288922837888845061849

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002837
$ws.Range("H5").Value = 0.008510999999999999
$ws.Range("I5").Value = 0.000071082477304929285628898317
$ws.Range("J5").Value = 0.000071082477304929285628898317
$ws.Range("M5").Value = 1.222391
$ws.Range("N5").Value = 3.667173
$ws.Range("O5").Value = 0.04717593986711188
$ws.Range("P5").Value = 0.04717593986711189
$ws.Range("Q5").Value = 0.003467923267
$ws.Range("R5").Value = 0.031211309403
$ws.Range("S5").Value = 0.000003353382674942688788770079
$ws.Range("T5").Value = 0.0000033533826749426900593195

$ws.Range("I6").Value = 0.3776915775490952
$ws.Range("J6").Value = 0.3776915775490952
$ws.Range("M6").Value = 1.815761
$ws.Range("N6").Value = 5.447283000000001
$ws.Range("O6").Value = 0.07007596730428067
$ws.Range("P6").Value = 0.07007596730428067
$ws.Range("Q6").Value = 27.37113236112667
$ws.Range("R6").Value = 246.34019125014
$ws.Range("S6").Value = 0.02646710263943259
$ws.Range("T6").Value = 0.02646710263943259

$ws.Range("I7").Value = 0.3776915775490952
$ws.Range("J7").Value = 0.3776915775490952
$ws.Range("O7").Value = 0.5079540516959071
$ws.Range("P7").Value = 0.5079540516959072
$ws.Range("S7").Value = 0.1918499671074818
$ws.Range("T7").Value = 0.1918499671074819

$ws.Range("I8").Value = 0.3776915775490952
$ws.Range("J8").Value = 0.3776915775490952
$ws.Range("M8").Value = 9.711409333333334
$ws.Range("N8").Value = 29.134228
$ws.Range("O8").Value = 0.3747940411327002
$ws.Range("P8").Value = 0.3747940411327002
$ws.Range("Q8").Value = 146.3916618298045
$ws.Range("R8").Value = 1317.52495646824
$ws.Range("S8").Value = 0.14155655265141
$ws.Range("T8").Value = 0.14155655265141

$ws.Range("I9").Value = 0.3776915775490952
$ws.Range("J9").Value = 0.3776915775490952
$ws.Range("M9").Value = 1.222391
$ws.Range("N9").Value = 3.667173
$ws.Range("O9").Value = 0.04717593986711188
$ws.Range("P9").Value = 0.04717593986711189
$ws.Range("Q9").Value = 18.42655826292667
$ws.Range("R9").Value = 165.83902436634
$ws.Range("S9").Value = 0.01781795515077074
$ws.Range("T9").Value = 0.01781795515077074

$ws.Range("G10").Value = 1.581618666666667
$ws.Range("H10").Value = 4.744856
$ws.Range("I10").Value = 0.03962825977384063
$ws.Range("J10").Value = 0.03962825977384063
$ws.Range("M10").Value = 1.815761
$ws.Range("N10").Value = 5.447283000000001
$ws.Range("O10").Value = 0.07007596730428067
$ws.Range("P10").Value = 0.07007596730428067
$ws.Range("Q10").Value = 2.871841491805334
$ws.Range("R10").Value = 25.846573426248
$ws.Range("S10").Value = 0.002776988636237197
$ws.Range("T10").Value = 0.002776988636237197

$ws.Range("G11").Value = 1.581618666666667
$ws.Range("H11").Value = 4.744856
$ws.Range("I11").Value = 0.03962825977384063
$ws.Range("J11").Value = 0.03962825977384063
$ws.Range("O11").Value = 0.5079540516959071
$ws.Range("P11").Value = 0.5079540516959072
$ws.Range("Q11").Value = 20.81688741101156
$ws.Range("R11").Value = 187.351986699104
$ws.Range("S11").Value = 0.02012933511378028
$ws.Range("T11").Value = 0.02012933511378029

$ws.Range("G12").Value = 1.581618666666667
$ws.Range("H12").Value = 4.744856
$ws.Range("I12").Value = 0.03962825977384063
$ws.Range("J12").Value = 0.03962825977384063
$ws.Range("M12").Value = 9.711409333333334
$ws.Range("N12").Value = 29.134228
$ws.Range("O12").Value = 0.3747940411327002
$ws.Range("P12").Value = 0.3747940411327002
$ws.Range("Q12").Value = 15.35974628124089
$ws.Range("R12").Value = 138.237716531168
$ws.Range("S12").Value = 0.01485243562369415
$ws.Range("T12").Value = 0.01485243562369415

$ws.Range("G13").Value = 1.581618666666667
$ws.Range("H13").Value = 4.744856
$ws.Range("I13").Value = 0.03962825977384063
$ws.Range("J13").Value = 0.03962825977384063
$ws.Range("M13").Value = 1.222391
$ws.Range("N13").Value = 3.667173
$ws.Range("O13").Value = 0.04717593986711188
$ws.Range("P13").Value = 0.04717593986711189
$ws.Range("Q13").Value = 1.933356423565333
$ws.Range("R13").Value = 17.400207812088
$ws.Range("S13").Value = 0.001869500400128994
$ws.Range("T13").Value = 0.001869500400128994

$ws.Range("G14").Value = 23.25273433333334
$ws.Range("H14").Value = 69.75820300000001
$ws.Range("I14").Value = 0.5826090801997593
$ws.Range("J14").Value = 0.5826090801997593
$ws.Range("M14").Value = 1.815761
$ws.Range("N14").Value = 5.447283000000001
$ws.Range("O14").Value = 0.07007596730428067
$ws.Range("P14").Value = 0.07007596730428067
$ws.Range("Q14").Value = 42.22140814582767
$ws.Range("R14").Value = 379.9926733124491
$ws.Range("S14").Value = 0.04082689485525537
$ws.Range("T14").Value = 0.04082689485525537

$ws.Range("G15").Value = 23.25273433333334
$ws.Range("H15").Value = 69.75820300000001
$ws.Range("I15").Value = 0.5826090801997593
$ws.Range("J15").Value = 0.5826090801997593
$ws.Range("O15").Value = 0.5079540516959071
$ws.Range("P15").Value = 0.5079540516959072
$ws.Range("Q15").Value = 306.0469396427391
$ws.Range("R15").Value = 2754.422456784652
$ws.Range("S15").Value = 0.2959386428422934
$ws.Range("T15").Value = 0.2959386428422935

$ws.Range("G16").Value = 23.25273433333334
$ws.Range("H16").Value = 69.75820300000001
$ws.Range("I16").Value = 0.5826090801997593
$ws.Range("J16").Value = 0.5826090801997593
$ws.Range("M16").Value = 9.711409333333334
$ws.Range("N16").Value = 29.134228
$ws.Range("O16").Value = 0.3747940411327002
$ws.Range("P16").Value = 0.3747940411327002
$ws.Range("Q16").Value = 225.8168212302538
$ws.Range("R16").Value = 2032.351391072284
$ws.Range("S16").Value = 0.2183584115686732
$ws.Range("T16").Value = 0.2183584115686732

$ws.Range("G17").Value = 23.25273433333334
$ws.Range("H17").Value = 69.75820300000001
$ws.Range("I17").Value = 0.5826090801997593
$ws.Range("J17").Value = 0.5826090801997593
$ws.Range("M17").Value = 1.222391
$ws.Range("N17").Value = 3.667173
$ws.Range("O17").Value = 0.04717593986711188
$ws.Range("P17").Value = 0.04717593986711189
$ws.Range("Q17").Value = 28.42393317445767
$ws.Range("R17").Value = 255.815398570119
$ws.Range("S17").Value = 0.02748513093353721
$ws.Range("T17").Value = 0.02748513093353721
